$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1 ("Problem Statement"): add szCs=28 (SizeBi=14) to the
#    paragraph mark run properties and the run itself, keeping the existing
#    sz=32 (Size=16) untouched.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of paragraph 4 ("Key objectives
#    include:") to the end of paragraph 1 ("Problem Statement"), right after
#    the run of text, before the paragraph mark.
#
#    Adding a zero-length bookmark exactly at "paragraph end - 1" confuses
#    this engine (it silently expands to the whole first paragraph), so we
#    temporarily insert a placeholder character after the target position,
#    anchor the bookmark right before it (now a "safe" offset), and then
#    remove the placeholder again. The bookmark stays put once the
#    placeholder text is deleted.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$endOfP1Text = $p1.Range.End - 1
$placeholder = $d.Range($endOfP1Text, $endOfP1Text)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endOfP1Text, $endOfP1Text)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($endOfP1Text, $endOfP1Text + 1)
$placeholderRange.Delete()

# ---------------------------------------------------------------------------
# 3) Paragraphs 2 and 3 (the two intro paragraphs): set sz=28, szCs=28
#    (Size=14, SizeBi=14) on the paragraph mark and every run.
# ---------------------------------------------------------------------------
foreach ($idx in 2, 3) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
}

# ---------------------------------------------------------------------------
# 4) Paragraph 4 ("Key objectives include:"): add sz=32, szCs=28 (Size=16,
#    SizeBi=14) alongside the existing bold formatting.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Font.Size = 16
$p4.Range.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 5) Paragraphs 5-9 (the four bullet items plus the closing paragraph): set
#    sz=28, szCs=28 (Size=14, SizeBi=14) on the paragraph mark and every run.
# ---------------------------------------------------------------------------
foreach ($idx in 5, 6, 7, 8, 9) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
}
